# Remove the parasitic (blank) row from the "2023-DRH-Annuel" sheet.
# The sheet had an extraneous empty row (row 2, between the header row and
# the first real data row) that carried only cell styles and no content.
# Deleting it shifts the real data rows up by one, matching the author's
# "retrait d'une ligne parasite" fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-DRH-Annuel")
$ws.Rows.Item(2).Delete()
